$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the contents of B14 with the new encrypted string value.
$ws.Range("B14").Value = "U2FsdGVkX1+35lx9lgFRMegpSQ9XGuPYSA3PLBTHm8poKFV+6PGHtQHWVwFvVCjGNgxlsz0x9DBtqPmu0CS8KS6U0XtfSLqoFG9LWH5Z1932EPUbRQFKFSo98Khq3ftd1Avxm2zK++wMotOAckVzyAdT0rWijbtUzqgzrui8P+M6ZLEEE1VP8oXOOPJ8Am0NYnAZR/1nizGEFuVJX6ty2uO/vb8J6VLn5dfoif0rBW+hOoa3ZcHrqzWWevrnhn25gpqkcHoXaqRO/Me65NRXiRj9BL7nox5mMH98wGGjTxD6J8eo5X61MC7qJp6ub85FkQ317MxDJjVtlDT7om0xW3+hdRjnTBZHF1nZgHY2etxc8un+go43Y4eNsRaeeo2fZFRVClYjHIVYeF+zkBVqo4dfc4rLY9B1uy6HPGI3ev/k4gD0OE0g6bk4260/O9Mxzr59bDztGemBrkMs3LbaRylZ6+DShcURTm1ldCfRoWgn1lyxdvn3ccpRqCabZseolfC9YxtIGpMudC4VtDtzxMlE0669rNqm9X5qzRAXMgG2KjfCEOtsnZJyblVyY7Cf1ipqcBL0tkuVFraFaj96Zt73DszwIJg1yrk6EIeGKP3uYB4B0Sxnk/m3LGMyM/t3pkvW2NNLLXmKn3WLQFCSQFW6PqMe48IbUSXoUlMjrmPOAQD605BgJCKT0PewkvkZWamHE6O7NYhOoRrDoQtFgEFULoYfpyHu0z9TF5THicEMj5G+JDLYlIEmmSzB+V2CWLV/kcJup8AzjZjzMY57Oq5FxH4DPFIWBx4h6xLJCM+KM3nPTXLRqrRxH7lIqWlsudZNpHCvhrD1VFGBhIOCPqxAaeiBLZu258kjJGHxvZGbEct5xvZkANLsISI8X2uaFVA5+OvA+E5OgbvBW71gkKwXGCHUjpjHIMegOG7fFQ9YfEbDIBcQlz6W8QwduKYWiYL2g3eg90RDd/6VEXMvCLvWTDzwZ2L6RhbqXlVq/vjFIdLr0RCEj5JNKsDlGjq42eFEZqwtKwIhq8nIY02mFR4CMP5ZF6wEoYz+2Ez96bGRQy+gAmZpNzBnxzLpRXTG"

# Move the active selection from B8 to B18, matching the saved view state.
$ws.Range("B18").Select()
